$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1788.5161
$ws.Range("I15").Value = 1788.5161
$ws.Range("K15").Value = 5365.5483
$ws.Range("M15").Value = -5196.5483
# Row 58
$ws.Range("H58").Value = 2412.6365
$ws.Range("I58").Value = 148.16667
$ws.Range("J58").Value = 5130
$ws.Range("K58").Value = 444.50001
$ws.Range("L58").Value = 15390
$ws.Range("M58").Value = -294.50001
$ws.Range("N58").Value = -15690
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 1874.75
$ws.Range("J11").Value = 1999.6666
$ws.Range("L11").Value = 1999.6666
$ws.Range("N11").Value = -2287.6666
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("N13").Value = 0
# Row 45
$ws.Range("H45").Value = 5998.875
$ws.Range("I45").Value = 5665.1665
$ws.Range("K45").Value = 5665.1665
$ws.Range("M45").Value = -5288.1665
# Row 74
$ws.Range("H74").Value = 13584.389
$ws.Range("I74").Value = 1798.625
$ws.Range("K74").Value = 1798.625
$ws.Range("M74").Value = -924.625
# Row 77
$ws.Range("H77").Value = 13584.389
$ws.Range("I77").Value = 1798.625
$ws.Range("K77").Value = 8993.125
$ws.Range("M77").Value = -4625.125
# Row 108
$ws.Range("H108").Value = 39620
$ws.Range("I108").Value = 39620
$ws.Range("K108").Value = 39620
$ws.Range("M108").Value = -35780
# Row 109
$ws.Range("H109").Value = 29000
$ws.Range("J109").Value = 29000
$ws.Range("L109").Value = 29000
$ws.Range("N109").Value = -31774
# Row 122
$ws.Range("H122").Value = 2035.8462
$ws.Range("I122").Value = 1810.5454
$ws.Range("K122").Value = 5431.6362
$ws.Range("M122").Value = -2981.6362
# Row 132
$ws.Range("H132").Value = 2274411.8
$ws.Range("J132").Value = 1084
$ws.Range("L132").Value = 3252
$ws.Range("N132").Value = -8312
# Row 133
$ws.Range("H133").Value = 69973.25
$ws.Range("J133").Value = 69973.25
$ws.Range("L133").Value = 69973.25
$ws.Range("N133").Value = -75033.25
# Row 135
$ws.Range("H135").Value = 79998.5
$ws.Range("J135").Value = 79998.5
$ws.Range("L135").Value = 79998.5
$ws.Range("N135").Value = -90138.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Range("H6").Value = 45737.25
$ws.Range("J6").Value = 45737.25
$ws.Range("L6").Value = 45737.25
$ws.Range("N6").Value = -45963.25
# Row 55
$ws.Range("H55").Value = 79500
$ws.Range("J55").Value = 79500
$ws.Range("L55").Value = 79500
$ws.Range("N55").Value = -80046
# Row 105
$ws.Range("H105").Value = 2434.0454
$ws.Range("I105").Value = 2082.2812
$ws.Range("K105").Value = 2082.2812
$ws.Range("M105").Value = -335.2811999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 63142.75
$ws.Range("I16").Value = 637.4
$ws.Range("K16").Value = 637.4
$ws.Range("M16").Value = -350.4
# Row 58
$ws.Range("H58").Value = 12033.637
$ws.Range("I58").Value = 2296.25
$ws.Range("J58").Value = 38000
$ws.Range("K58").Value = 2296.25
$ws.Range("L58").Value = 38000
$ws.Range("M58").Value = -2093.25
$ws.Range("N58").Value = -38406
# Row 92
$ws.Range("H92").Value = 50300.5
$ws.Range("J92").Value = 50300.5
$ws.Range("L92").Value = 50300.5
$ws.Range("N92").Value = -55292.5
# Row 113
$ws.Range("H113").Value = 63142.75
$ws.Range("I113").Value = 637.4
$ws.Range("K113").Value = 637.4
$ws.Range("M113").Value = 1532.6
# Row 114
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678
# Row 119
$ws.Range("H119").Value = 19845.5
$ws.Range("I119").Value = 19691
$ws.Range("J119").Value = 20000
$ws.Range("K119").Value = 19691
$ws.Range("L119").Value = 20000
$ws.Range("M119").Value = -14853
$ws.Range("N119").Value = -29676
# Row 134
$ws.Range("H134").Value = 2039.125
$ws.Range("I134").Value = 1564.7142
$ws.Range("J134").Value = 5360
$ws.Range("K134").Value = 4694.142599999999
$ws.Range("L134").Value = 16080
$ws.Range("M134").Value = -2159.142599999999
$ws.Range("N134").Value = -21150
# Row 136
$ws.Range("H136").Value = 12033.637
$ws.Range("I136").Value = 2296.25
$ws.Range("J136").Value = 38000
$ws.Range("K136").Value = 6888.75
$ws.Range("L136").Value = 114000
$ws.Range("M136").Value = -4338.75
$ws.Range("N136").Value = -119100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 5793.9287
$ws.Range("I109").Value = 1019.1667
$ws.Range("K109").Value = 3057.5001
$ws.Range("M109").Value = -2017.5001
# Row 131
$ws.Range("H131").Value = 9142.857
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
# Row 137
$ws.Range("H137").Value = 9513.799999999999
$ws.Range("J137").Value = 12589.8
$ws.Range("L137").Value = 37769.39999999999
$ws.Range("N137").Value = -47969.39999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 66.052635
$ws.Range("I2").Value = 41.5
$ws.Range("K2").Value = 41.5
$ws.Range("M2").Value = 71.5
# Row 7
$ws.Range("H7").Value = 1002500
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224
# Row 8
$ws.Range("H8").Value = 1002500
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5278
# Row 11
$ws.Range("H11").Value = 4282892
$ws.Range("I11").Value = 6722971.5
$ws.Range("J11").Value = 12752.75
$ws.Range("K11").Value = 6722971.5
$ws.Range("L11").Value = 12752.75
$ws.Range("M11").Value = -6722832.5
$ws.Range("N11").Value = -13030.75
# Row 95
$ws.Range("H95").Value = 39000
$ws.Range("J95").Value = 39000
$ws.Range("L95").Value = 39000
$ws.Range("N95").Value = -44492

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 4203.2856
$ws.Range("I93").Value = 2191.5
$ws.Range("K93").Value = 2191.5
$ws.Range("M93").Value = -943.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 9000
$ws.Range("J10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("N10").Value = -9338
# Row 24
$ws.Range("H24").Value = 504504.5
$ws.Range("I24").Value = 504504.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 504504.5
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0
$ws.Range("M24").Value = -504274.5
# Row 96
$ws.Range("H96").Value = 2877.182
$ws.Range("I96").Value = 2187.25
$ws.Range("K96").Value = 2187.25
$ws.Range("M96").Value = -814.25
# Row 103
$ws.Range("H103").Value = 700
$ws.Range("J103").Value = 700
$ws.Range("L103").Value = 700
$ws.Range("N103").Value = -3044
# Row 123
$ws.Range("H123").Value = 47726.637
$ws.Range("J123").Value = 47726.637
$ws.Range("L123").Value = 47726.637
$ws.Range("N123").Value = -57526.637
# Row 136
$ws.Range("H136").Value = 171272.16
$ws.Range("I136").Value = 224753.62
$ws.Range("J136").Value = 1914.1666
$ws.Range("K136").Value = 674260.86
$ws.Range("L136").Value = 5742.4998
$ws.Range("M136").Value = -671710.86
$ws.Range("N136").Value = -10842.4998
